$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing the existing rows 6-8 down to 7-9.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new observation record (Ullticka / Phellinidium ferrugineofuscum).
$ws.Cells.Item(6, 1).Value() = 112044152
$ws.Cells.Item(6, 2).Value() = 89405
$ws.Cells.Item(6, 3).Value() = "Ovaliderad"
$ws.Cells.Item(6, 4).Value() = "NT"
$ws.Cells.Item(6, 5).Value() = 1202
$ws.Cells.Item(6, 6).Value() = "Ullticka"
$ws.Cells.Item(6, 7).Value() = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(6, 8).Value() = "(P.Karst.) Fiasson & Niemelä"

$ws.Cells.Item(6, 16).Value() = "Kyrkberget, Dlr"
$ws.Cells.Item(6, 17).Value() = 554745.9054445035
$ws.Cells.Item(6, 18).Value() = 6697510.461741267
$ws.Cells.Item(6, 19).Value() = 15
$ws.Cells.Item(6, 20).Value() = "Dalarna"
$ws.Cells.Item(6, 21).Value() = "Hedemora"
$ws.Cells.Item(6, 22).Value() = "Dalarna"
$ws.Cells.Item(6, 23).Value() = "Husby"

# Dates are stored as literal text in this sheet, so quote-prefix them to stop
# them being auto-coerced into date serial numbers.
$ws.Cells.Item(6, 25).Value() = "'2023-09-12"
$ws.Cells.Item(6, 26).Value() = "00:00"
$ws.Cells.Item(6, 27).Value() = "'2023-09-12"
$ws.Cells.Item(6, 28).Value() = "00:00"

$ws.Cells.Item(6, 30).Value() = $False
$ws.Cells.Item(6, 31).Value() = $False
$ws.Cells.Item(6, 33).Value() = $False

$ws.Cells.Item(6, 49).Value() = "Philipp Weiss"
$ws.Cells.Item(6, 50).Value() = "Philipp Weiss"
